# Schematic updates. First export to PCB
# Added 9V line, OLED switch circuit. Added OLED footprint.
#
# Adds the new "78L09-150" 9V LDO part as a new row at the bottom of the
# BOM table (Table1), wires up its hyperlink to the LCSC product page, and
# moves the active selection the way the author's session ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BOM rows live in an Excel Table (ListObject) - grow it by one row so
# the table range + autofilter + dimension all expand together, just like
# clicking "Insert Table Row" at the bottom of Table1 would.
$tbl = $ws.ListObjects.Item(1)
$newListRow = $tbl.ListRows.Add()
$newRow = $newListRow.Range.Row

# Fill columns in the same order the workbook's shared-string table shows
# them having been entered (Link text first, then Description, then Part),
# so new unique strings land in the expected slots.
$linkText = "78L09-150 UMW(Youtai Semiconductor Co., Ltd.) | C347271 - LCSC Electronics"
$ws.Cells.Item($newRow, 6).Value = $linkText
$ws.Cells.Item($newRow, 2).Value = "9V LDO"
$ws.Cells.Item($newRow, 1).Value = "78L09-150"
$ws.Cells.Item($newRow, 3).Value = 0.43
$ws.Cells.Item($newRow, 4).Value = 1
$ws.Cells.Item($newRow, 5).Formula = "=Table1[[#This Row],[Price]]*Table1[[#This Row],[Qty]]"

# Hyperlink the Link cell to the LCSC product page, keeping the cell's
# text as the displayed text.
$linkCell = $ws.Cells.Item($newRow, 6)
$lcscUrl = "https://www.lcsc.com/product-detail/Linear-Voltage-Regulators-LDO_UMW-Youtai-Semiconductor-Co-Ltd-78L09-150_C347271.html"
$hlink = $ws.Hyperlinks.Add($linkCell, $lcscUrl, [Type]::Missing, [Type]::Missing, $linkText)
# Match the existing rows' hyperlink bookkeeping, where the stored
# "display" (used when the target can't be reached) echoes the URL.
$hlink.TextToDisplay = $lcscUrl
$linkCell.Value = $linkText
# Re-apply the workbook's Hyperlink cell style explicitly (Hyperlinks.Add
# already visually formats the cell, but this keeps it on the workbook's
# shared "Hyperlink" style instead of a freshly minted duplicate one).
$linkCell.Style = "Hyperlink"

# Match the ending selection left behind in the saved session.
$ws.Range("A23").Select() | Out-Null
